# Applies updated cryptocurrency price/volume figures to Sheet1
# per commit "Updated symbol list on Thu Feb  2 15:39:04 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.02%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.78%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.696"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "11.20%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08093"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.86%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.88%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.686"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.61%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.963"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.13%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9444"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.21%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1281"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "16.20%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1989"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.00%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09226"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.12%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03419"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09605"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.32%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001321"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.23%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006064"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.26%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.372"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.54%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3507"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.58%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.598"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9.26%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04443"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.25%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001253"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.95%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001191"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-15.02%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003995"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "37.63%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02528"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "17.35%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05211"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.88%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007358"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.79%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1430"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.64%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009028"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002192"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.69%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01002"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "25.68%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006697"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.69%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002875"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-12.76%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001803"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "24.80%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
